# Table4_PublishedResults_method.docx — correct ordering of factor
# variables: the page's section is switched from portrait to landscape
# orientation so the (wide) results table fits the page.
$d = $word.ActiveDocument

# wdOrientLandscape = 1 (wdOrientPortrait = 0). Word automatically swaps
# PageSetup.PageWidth/PageHeight to match the new orientation, mirroring
# the pgSz w:w/h swap + orient="landscape" seen in the target XML.
$d.PageSetup.Orientation = 1
